$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "65.250.44"
$ws.Range("E2").Value2 = "  -0.32%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.541.90"
$ws.Range("E3").Value2 = "  +3.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  -0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "597.12"
$ws.Range("E5").Value2 = "  +1.70%  "

# Row 6 - Solana
Set-TextValue "D6" "138.79"
$ws.Range("E6").Value2 = "  +0.83%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.541.11"
$ws.Range("E7").Value2 = "  +3.36%  "

# Row 8 - USDC
$ws.Range("E8").Value2 = "  +0.07%  "

# Row 9 - XRP
Set-TextValue "D9" "0.494"
$ws.Range("E9").Value2 = "  -1.03%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value2 = "  +3.12%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.89"
$ws.Range("E11").Value2 = "  -5.17%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.387"
$ws.Range("E12").Value2 = "  +3.08%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.143.23"
$ws.Range("E13").Value2 = "  +3.36%  "

# Row 14 - ShibaInu
Set-TextValue "D14" "0.0000185"
$ws.Range("E14").Value2 = "  +2.78%  "

# Row 15 - Avalanche
$ws.Range("E15").Value2 = "  +3.61%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.544.63"
$ws.Range("E16").Value2 = "  +3.31%  "

# Row 17 - TRON
$ws.Range("E17").Value2 = "  +1.51%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "65.254.44"
$ws.Range("E18").Value2 = "  -0.29%  "

# Row 19 - Uniswap
Set-TextValue "D19" "10.25"
$ws.Range("E19").Value2 = "  +4.72%  "

# Row 20 - Polkadot
Set-TextValue "D20" "5.95"
$ws.Range("E20").Value2 = "  +1.26%  "

# Row 21 - Chainlink
Set-TextValue "D21" "14.36"

# Row 22 - BitcoinCash
Set-TextValue "D22" "394.11"
$ws.Range("E22").Value2 = "  +0.75%  "

# Row 23 - Polygon
$ws.Range("E23").Value2 = "  +3.22%  "

# Row 24 - WrappedeETH
Set-TextValue "D24" "3.683.82"
$ws.Range("E24").Value2 = "  +3.34%  "

# Row 25 - Litecoin
Set-TextValue "D25" "73.83"
$ws.Range("E25").Value2 = "  +0.89%  "

# Row 26 - Dai
$ws.Range("E26").Value2 = "  -0.04%  "

# Row 27 - PEPE
$ws.Range("E27").Value2 = "  +8.74%  "

# Row 28 - RenderToken
Set-TextValue "D28" "7.81"
$ws.Range("E28").Value2 = "  +9.10%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "D29" "0.998"
$ws.Range("E29").Value2 = "  +0.20%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value2 = "  +2.50%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "8.15"
$ws.Range("E31").Value2 = "  -0.45%  "

# Row 32 - RenzoRestakedETH
Set-TextValue "D32" "3.553.83"
$ws.Range("E32").Value2 = "  +3.48%  "

# Row 33 - USDe
$ws.Range("E33").Value2 = "  +0.03%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "23.86"
$ws.Range("E34").Value2 = "  +3.87%  "

# Row 35 - Kaspa
$ws.Range("E35").Value2 = "  +0.95%  "

# Row 36 - Fetch.AI
Set-TextValue "D36" "1.29"
$ws.Range("E36").Value2 = "  +11.46%  "

# Row 37 - Aptos
Set-TextValue "D37" "6.97"

# Row 38 - Monero
Set-TextValue "D38" "169.23"
$ws.Range("E38").Value2 = "  -2.17%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value2 = "  +5.59%  "

# Row 40 - NEARProtocol
$ws.Range("E40").Value2 = "  +5.44%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0803"
$ws.Range("E41").Value2 = "  +4.89%  "

# Row 42 - Mantle
$ws.Range("E42").Value2 = "  +0.90%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "26.39"
$ws.Range("E43").Value2 = "  +18.06%  "

# Row 44 - OKB
Set-TextValue "D44" "42.76"
$ws.Range("E44").Value2 = "  -1.81%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value2 = "  -0.06%  "

# Row 46 - Filecoin
Set-TextValue "D46" "4.43"
$ws.Range("E46").Value2 = "  +0.09%  "

# Row 47 - ONDO
$ws.Range("E47").Value2 = "  +7.88%  "

# Row 48 - Stacks
$ws.Range("E48").Value2 = "  +3.80%  "

# Row 49 - Cosmos
$ws.Range("E49").Value2 = "  +4.17%  "

# Row 50 - Maker
Set-TextValue "D50" "2.396.78"
$ws.Range("E50").Value2 = "  +9.28%  "

# Row 51 - Bittensor
Set-TextValue "D51" "306.91"
$ws.Range("E51").Value2 = "  +6.99%  "
